# Insert a new weekly price record for "Perejil" at row 200 in the
# "Feria Lagunitas de Puerto Montt" sheet. Excel's native Rows().Insert()
# shifts the existing rows 200-261 down to 201-262 (carrying their values
# and formatting with them), matching the target diff exactly. We then
# populate the freshly inserted, empty row 200 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(200).Insert()

$ws.Range("A200").Value = 4
$ws.Range("B200").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C200").Value = "Los Lagos"
$ws.Range("D200").Value = 44736
$ws.Range("E200").Value = 10
$ws.Range("F200").Value = 100112044
$ws.Range("G200").Value = "Perejil"
$ws.Range("H200").Value = "Sin especificar"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 180
$ws.Range("K200").Value = 6000
$ws.Range("L200").Value = 6000
$ws.Range("M200").Value = 6000
$ws.Range("N200").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O200").Value = "Región de La Araucanía"
$ws.Range("P200").Value = 3000
$ws.Range("Q200").Value = 2
$ws.Range("R200").Value = "Hortaliza"
